# Insert a new record row at row 155 ("Femacal de La Calera" / Berenjena
# price sheet). This pushes the previous rows 155-194 down to 156-195,
# and fills the newly inserted row with a new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 155, shifting existing
# data (old rows 155-194) down to rows 156-195.
$ws.Rows.Item(155).Insert()

# Populate the new row 155 with the new observation. Most columns carry
# the same constant values used throughout this subset (Mercado,
# Región, Categoría, etc.); only the date, volume, prices and resulting
# $/Kg differ.
$ws.Cells.Item(155, 1).Value  = 3
$ws.Cells.Item(155, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(155, 3).Value  = "Coquimbo"
$ws.Cells.Item(155, 4).Value  = 44543
$ws.Cells.Item(155, 5).Value  = 5
$ws.Cells.Item(155, 6).Value  = 100112001
$ws.Cells.Item(155, 7).Value  = "Berenjena"
$ws.Cells.Item(155, 8).Value  = "Sin especificar"
$ws.Cells.Item(155, 9).Value  = "Primera"
$ws.Cells.Item(155, 10).Value = 105
$ws.Cells.Item(155, 11).Value = 8500
$ws.Cells.Item(155, 12).Value = 9000
$ws.Cells.Item(155, 13).Value = 8738
$ws.Cells.Item(155, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(155, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(155, 16).Value = 146
$ws.Cells.Item(155, 17).Value = 60
$ws.Cells.Item(155, 18).Value = "Hortaliza"
